# ajout du datapath de l'instruction jal dans le schema et des CU codes dans la table de verite
#
# This script reproduces, on the "table de verite" worksheet, the switch of the
# RegDst (col D) and MemToReg (col H) columns from raw 0/1 numbers to the same
# "00"/"01"/"X" textual code convention already used for the other control-unit
# columns, applies a new orange font to the newly-recoded cells, and appends a
# new row describing the "jal" instruction (plus a couple of stray formatted
# cells below it, as Excel leaves behind when a row is built then trimmed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Orange font color used for the newly recoded "CU codes" (RGB FFF78E25 -> BGR long)
$orange = 2461431

# ---------------------------------------------------------------------------
# Step 1: cells that only need the new orange font (value already "X" text) -
# touch these FIRST so the non-quotePrefix orange style is allocated before
# the quotePrefix one.
# ---------------------------------------------------------------------------
foreach ($addr in @("D4", "H4", "D5", "H5", "D6", "H6")) {
    $c = $ws.Range($addr)
    $c.Font.Color = $orange
}

# ---------------------------------------------------------------------------
# Step 2: cells whose numeric 0/1 becomes a quoted textual "00"/"01" code,
# recolored orange.
# ---------------------------------------------------------------------------
$c = $ws.Range("D2"); $c.Value = "'01"; $c.Font.Color = $orange
$c = $ws.Range("H2"); $c.Value = "'00"; $c.Font.Color = $orange

$c = $ws.Range("D3"); $c.Value = "'00"; $c.Font.Color = $orange
$c = $ws.Range("H3"); $c.Value = "'01"; $c.Font.Color = $orange

$c = $ws.Range("D7"); $c.Value = "'00"; $c.Font.Color = $orange
$c = $ws.Range("H7"); $c.Value = "'00"; $c.Font.Color = $orange

# Row 8 keeps its existing blue font, only the value becomes a quoted "00".
$c = $ws.Range("D8"); $c.Value = "'00"
$c = $ws.Range("H8"); $c.Value = "'00"

# ---------------------------------------------------------------------------
# Step 3: new row 10 - the "jal" instruction line.
# ---------------------------------------------------------------------------
$c = $ws.Range("A10"); $c.Value = "jal";     $c.Font.Color = $orange
$c = $ws.Range("B10"); $c.Value = "'000011"; $c.Font.Color = $orange
$c = $ws.Range("C10"); $c.Value = 1;         $c.Font.Color = $orange
$c = $ws.Range("D10"); $c.Value = "'10";     $c.Font.Color = $orange
$c = $ws.Range("E10"); $c.Value = "X";       $c.Font.Color = $orange
$c = $ws.Range("F10"); $c.Value = "X";       $c.Font.Color = $orange
$c = $ws.Range("G10"); $c.Value = 0;         $c.Font.Color = $orange
$c = $ws.Range("H10"); $c.Value = "'10";     $c.Font.Color = $orange
$c = $ws.Range("I10"); $c.Value = "XX";      $c.Font.Color = $orange
$c = $ws.Range("J10"); $c.Value = "'01";     $c.Font.Color = $orange

# Trailing formatted-but-empty cells left over in the saved sheet.
$ws.Range("K10").Font.Color = $orange
$ws.Range("G14").Font.Color = $orange

# ---------------------------------------------------------------------------
# Step 4: page setup + selection, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

[void]$ws.Range("E17").Select()
